$wb = $excel.ActiveWorkbook

# --- Sheet "展览": shift rows 4-14 content up into rows 2-12 (cols B-I), zero col F, then drop trailing rows 13-14 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("B2").Value = '2024-07-12'
$ws1.Range("C2").Value = '南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展'
$ws1.Range("D2").Value = '民族大道106号 南宁国际会展中心'
$ws1.Range("E2").Value = '2024.07.12 09:30-07.14 17:00'
$ws1.Range("F2").Value = 0
$ws1.Range("G2").Value = 50
$ws1.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=87182'
$ws1.Range("I2").Value = '//i0.hdslb.com/bfs/openplatform/202406/CsYbpZmU1719311879090.jpeg'

$ws1.Range("B3").Value = '2024-07-13'
$ws1.Range("C3").Value = '南宁·0713国乙ONLY'
$ws1.Range("D3").Value = '亭洪路45号 水明漾宴会中心'
$ws1.Range("E3").Value = '2024.07.13 09:30-07.13 21:00'
$ws1.Range("F3").Value = 0
$ws1.Range("G3").Value = 68
$ws1.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=86378'
$ws1.Range("I3").Value = '//i1.hdslb.com/bfs/openplatform/202405/ZDBCv2of1716659486569.jpeg'

$ws1.Range("B4").Value = '2024-07-14'
$ws1.Range("C4").Value = '广西·首届明日方舟only展 - 花庭圣梦'
$ws1.Range("D4").Value = '明秀东路157号 利泰国际大酒店'
$ws1.Range("E4").Value = '2024.07.14 09:00-07.14 18:00'
$ws1.Range("F4").Value = 0
$ws1.Range("G4").Value = 69
$ws1.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=85852'
$ws1.Range("I4").Value = '//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg'

$ws1.Range("B5").Value = '2024-07-20'
$ws1.Range("C5").Value = '南宁·AB动漫游戏嘉年华'
$ws1.Range("D5").Value = '五象大道西段669号 广西体育中心体育馆'
$ws1.Range("E5").Value = '2024.07.20 09:30-07.21 17:00'
$ws1.Range("F5").Value = 0
$ws1.Range("G5").Value = 60
$ws1.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=84862'
$ws1.Range("I5").Value = '//i1.hdslb.com/bfs/openplatform/202407/R7iP9Iio1720170437964.jpeg'

$ws1.Range("B6").Value = '2024-07-20'
$ws1.Range("C6").Value = '横州·第二届海棠动漫游戏嘉年华'
$ws1.Range("D6").Value = '茉莉花大道 横州国际大酒店'
$ws1.Range("E6").Value = '2024.07.20 09:30-07.20 17:00'
$ws1.Range("F6").Value = 0
$ws1.Range("G6").Value = 30
$ws1.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=84799'
$ws1.Range("I6").Value = '//i2.hdslb.com/bfs/openplatform/202404/r50S2ttT1713869164413.jpeg'

$ws1.Range("B7").Value = '2024-07-27'
$ws1.Range("C7").Value = '南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）'
$ws1.Range("D7").Value = '民族大道106号 南宁国际会展中心'
$ws1.Range("E7").Value = '2024.07.27 09:30-07.28 17:30'
$ws1.Range("F7").Value = 0
$ws1.Range("G7").Value = 55
$ws1.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=85264'
$ws1.Range("I7").Value = '//i1.hdslb.com/bfs/openplatform/202406/JxFed5iv1718622152091.jpeg'

$ws1.Range("B8").Value = '2024-08-03'
$ws1.Range("C8").Value = '南宁·火影忍者only'
$ws1.Range("D8").Value = '厢竹大道65号 桔子酒店'
$ws1.Range("E8").Value = '2024.08.03 10:00-08.03 17:00'
$ws1.Range("F8").Value = 0
$ws1.Range("G8").Value = 68
$ws1.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=86994'
$ws1.Range("I8").Value = '//i0.hdslb.com/bfs/openplatform/202406/h1tXE9t11717523356034.jpeg'

$ws1.Range("B9").Value = '2024-08-03'
$ws1.Range("C9").Value = '南宁·蔚蓝档案only'
$ws1.Range("D9").Value = '亭洪路45号 百益上河城'
$ws1.Range("E9").Value = '2024.08.03 09:00-08.03 17:00'
$ws1.Range("F9").Value = 0
$ws1.Range("G9").Value = 68
$ws1.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=85370'
$ws1.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202405/sBxi2Mx61715247424836.jpeg'

$ws1.Range("B10").Value = '2024-08-10'
$ws1.Range("C10").Value = '南宁·国乙only'
$ws1.Range("D10").Value = '新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店'
$ws1.Range("E10").Value = '2024.08.10 10:00-08.10 17:00'
$ws1.Range("F10").Value = 0
$ws1.Range("G10").Value = 40
$ws1.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=88227'
$ws1.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202406/3cFX9LLQ1719482186347.jpeg'

$ws1.Range("B11").Value = '2024-08-24'
$ws1.Range("C11").Value = '南宁·第二届北极光动漫展'
$ws1.Range("D11").Value = '民族大道106号 南宁国际会展中心'
$ws1.Range("E11").Value = '2024.08.24 09:00-08.25 17:00'
$ws1.Range("F11").Value = 0
$ws1.Range("G11").Value = 65
$ws1.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=88276'
$ws1.Range("I11").Value = '//i1.hdslb.com/bfs/openplatform/202406/mTEwC1GY1717576221099.jpeg'

$ws1.Range("B12").Value = '2024-11-02'
$ws1.Range("C12").Value = '南宁·万圣漫控嘉年华10'
$ws1.Range("D12").Value = '亭洪路45号 百益上河城'
$ws1.Range("E12").Value = '2024.11.02 11:00-11.03 22:00'
$ws1.Range("F12").Value = 0
$ws1.Range("G12").Value = 50
$ws1.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=87820'
$ws1.Range("I12").Value = '//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg'

# remove now-duplicate trailing rows 13 and 14
$ws1.Rows("13:14").Delete()

# --- Sheet "演出": zero out col F (rows 2-6) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 0
$ws2.Range("F3").Value = 0
$ws2.Range("F4").Value = 0
$ws2.Range("F5").Value = 0
$ws2.Range("F6").Value = 0

# --- Sheet "全部类型": zero out col F (rows 2-19); G2/G3 become text "已停售" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("G2").Value = '已停售'
$ws4.Range("F3").Value = 0
$ws4.Range("G3").Value = '已停售'
$ws4.Range("F4").Value = 0
$ws4.Range("F5").Value = 0
$ws4.Range("F6").Value = 0
$ws4.Range("F7").Value = 0
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 0
$ws4.Range("F12").Value = 0
$ws4.Range("F13").Value = 0
$ws4.Range("F14").Value = 0
$ws4.Range("F15").Value = 0
$ws4.Range("F16").Value = 0
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 0
$ws4.Range("F19").Value = 0
